$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.55
$ws.Range("I2").Value = 3
$ws.Range("X2").Value = 12
$ws.Range("AH2").Value = 8.5
$ws.Range("AI2").Value = 13
$ws.Range("AJ2").Value = 11
$ws.Range("AK2").Value = 29
$ws.Range("AN2").Value = 4.5
